# Adapt the workbook to the renamed "mydictionary" repository:
#   - insert a new "Note" column between "Define" and "QC"
#   - keep the autofilter / filter-database range in sync with the new extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts the existing "QC"/"QT" columns
# (D,E) one place to the right, to E and F.
$ws.Columns.Item(4).Insert()

# Match the new column's width to its left neighbour ("Define", column C)
# so the new column reads comfortably.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Give the new column its header text.
$ws.Range("D1").Value = "Note"

# Re-establish the AutoFilter over the widened header row (A1:F1).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F1").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the
# new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='bing-dictionary'!`$A`$1:`$F`$1"
    }
}
